$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 14, shifting existing rows 14-69 down to 15-70.
$ws.Rows("14").Insert()

# Populate the new row 14 with the new record's data.
$ws.Cells.Item(14, 1).Value = 5
$ws.Cells.Item(14, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(14, 3).Value = "Maule"
$ws.Cells.Item(14, 4).Value = 44525
$ws.Cells.Item(14, 5).Value = 7
$ws.Cells.Item(14, 6).Value = 100112022
$ws.Cells.Item(14, 7).Value = "Arveja Verde"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 400
$ws.Cells.Item(14, 11).Value = 14000
$ws.Cells.Item(14, 12).Value = 14000
$ws.Cells.Item(14, 13).Value = 14000
$ws.Cells.Item(14, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(14, 15).Value = "Región del Maule"
$ws.Cells.Item(14, 16).Value = 560
$ws.Cells.Item(14, 17).Value = 25
$ws.Cells.Item(14, 18).Value = "Hortaliza"
